# Fluent wait, update property files, dynamic xpath
#
# The "Password" row's object-repository value was obfuscated/updated
# (pwd -> p@@d) and a new "password" label was added alongside it in
# column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 previously held the literal "pwd" placeholder; update it.
$ws.Range("C2").Value = "p@@d"

# New column D, row 2: add the "password" label.
$ws.Range("D2").Value = "password"

# Leave the selection on the newly-edited cell, like Excel would after
# typing the value and pressing Enter/Tab.
$ws.Range("D2").Select()
